$wb = $excel.ActiveWorkbook

# Worksheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 839.7692
$ws.Range("I33").Value = 966.1429000000001
$ws.Range("K33").Value = 966.1429000000001
$ws.Range("M33").Value = -737.1429000000001
$ws.Range("H96").Value = 652.0769
$ws.Range("I96").Value = 286.33334
$ws.Range("J96").Value = 1475
$ws.Range("K96").Value = 859.0000200000001
$ws.Range("L96").Value = 4425
$ws.Range("M96").Value = 513.9999799999999
$ws.Range("N96").Value = -7171
$ws.Range("H112").Value = 3448.1428
$ws.Range("I112").Value = 1714
$ws.Range("J112").Value = 4748.75
$ws.Range("K112").Value = 5142
$ws.Range("L112").Value = 14246.25
$ws.Range("M112").Value = -4034
$ws.Range("N112").Value = -16462.25
$ws.Range("H138").Value = 10957.612
$ws.Range("J138").Value = 11045.406
$ws.Range("L138").Value = 33136.218
$ws.Range("N138").Value = -43416.218

# Worksheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 6127.5293
$ws.Range("I61").Value = 5478.5386
$ws.Range("K61").Value = 5478.5386
$ws.Range("M61").Value = -5266.5386
$ws.Range("H110").Value = 1678.5
$ws.Range("I110").Value = 1534.1351
$ws.Range("J110").Value = 2441.5715
$ws.Range("K110").Value = 1534.1351
$ws.Range("L110").Value = 2441.5715
$ws.Range("M110").Value = 510.8649
$ws.Range("N110").Value = -6531.5715
$ws.Range("H136").Value = 6127.5293
$ws.Range("I136").Value = 5478.5386
$ws.Range("K136").Value = 16435.6158
$ws.Range("M136").Value = -13885.6158

# Worksheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1917.6522
$ws.Range("I20").Value = 1690.35
$ws.Range("K20").Value = 1690.35
$ws.Range("M20").Value = -1443.35
$ws.Range("H134").Value = 3050
$ws.Range("I134").Value = 2260.2
$ws.Range("K134").Value = 6780.599999999999
$ws.Range("M134").Value = -4245.599999999999

# Worksheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 6915
$ws.Range("J31").Value = 10566.667
$ws.Range("L31").Value = 10566.667
$ws.Range("N31").Value = -11156.667
$ws.Range("H34").Value = 6915
$ws.Range("J34").Value = 10566.667
$ws.Range("L34").Value = 10566.667
$ws.Range("N34").Value = -10970.667
$ws.Range("H88").Value = 27199.8
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 27199.8
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 27199.8
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -28011.8
$ws.Range("H91").Value = 27199.8
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 27199.8
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 27199.8
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -30007.8
$ws.Range("H94").Value = 673.6
$ws.Range("I94").Value = 622
$ws.Range("J94").Value = 880
$ws.Range("K94").Value = 622
$ws.Range("L94").Value = 880
$ws.Range("M94").Value = -171
$ws.Range("N94").Value = -1782
$ws.Range("H132").Value = 3423.9092
$ws.Range("I132").Value = 2995.875
$ws.Range("J132").Value = 4565.3335
$ws.Range("K132").Value = 8987.625
$ws.Range("L132").Value = 13696.0005
$ws.Range("M132").Value = -6457.625
$ws.Range("N132").Value = -18756.0005
$ws.Range("H133").Value = 86876.06
$ws.Range("J133").Value = 87930.81
$ws.Range("L133").Value = 87930.81
$ws.Range("N133").Value = -92990.81
$ws.Range("H134").Value = 4999.5
$ws.Range("I134").Value = 4999.5
$ws.Range("K134").Value = 14998.5
$ws.Range("M134").Value = -12463.5
$ws.Range("H141").Value = 215034.83
$ws.Range("I141").Value = 50000
$ws.Range("K141").Value = 50000
$ws.Range("M141").Value = -44820

# Worksheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H42").Value = 4324.3335
$ws.Range("J42").Value = 4986.5
$ws.Range("L42").Value = 14959.5
$ws.Range("N42").Value = -16027.5
$ws.Range("H98").Value = 997
$ws.Range("J98").Value = 1109
$ws.Range("L98").Value = 3327
$ws.Range("N98").Value = -6323
$ws.Range("H131").Value = 3711.7778
$ws.Range("J131").Value = 4401.5
$ws.Range("L131").Value = 13204.5
$ws.Range("N131").Value = -23284.5
$ws.Range("H139").Value = 4782.625
$ws.Range("J139").Value = 5247.6875
$ws.Range("L139").Value = 15743.0625
$ws.Range("N139").Value = -26023.0625

# Worksheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 6364.4287
$ws.Range("I70").Value = 6370
$ws.Range("K70").Value = 6370
$ws.Range("M70").Value = -6100
$ws.Range("H73").Value = 6364.4287
$ws.Range("I73").Value = 6370
$ws.Range("K73").Value = 6370
$ws.Range("M73").Value = -5434
$ws.Range("H80").Value = 7341.893
$ws.Range("I80").Value = 5944.727
$ws.Range("J80").Value = 8245.941000000001
$ws.Range("K80").Value = 5944.727
$ws.Range("L80").Value = 8245.941000000001
$ws.Range("M80").Value = -4946.727
$ws.Range("N80").Value = -10241.941
$ws.Range("H83").Value = 7341.893
$ws.Range("I83").Value = 5944.727
$ws.Range("J83").Value = 8245.941000000001
$ws.Range("K83").Value = 29723.635
$ws.Range("L83").Value = 41229.705
$ws.Range("M83").Value = -24731.635
$ws.Range("N83").Value = -51213.705
$ws.Range("H123").Value = 54391.168
$ws.Range("J123").Value = 54391.168
$ws.Range("L123").Value = 54391.168
$ws.Range("N123").Value = -59291.168
$ws.Range("H126").Value = 6942.0713
$ws.Range("I126").Value = 6547.5
$ws.Range("K126").Value = 19642.5
$ws.Range("M126").Value = -17172.5

# Worksheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 14672.3
$ws.Range("I7").Value = 28800.75
$ws.Range("K7").Value = 28800.75
$ws.Range("M7").Value = -28688.75
$ws.Range("H21").Value = 14000
$ws.Range("I21").Value = 14000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 14000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -13826
$ws.Range("N21").ClearContents()
$ws.Range("H40").Value = 15184.538
$ws.Range("I40").Value = 14665.444
$ws.Range("K40").Value = 14665.444
$ws.Range("M40").Value = -14529.444
$ws.Range("H82").Value = 1808.3529
$ws.Range("I82").Value = 1703
$ws.Range("J82").Value = 2001.5
$ws.Range("K82").Value = 1703
$ws.Range("L82").Value = 2001.5
$ws.Range("M82").Value = -1342
$ws.Range("N82").Value = -2723.5
$ws.Range("H85").Value = 1808.3529
$ws.Range("I85").Value = 1703
$ws.Range("J85").Value = 2001.5
$ws.Range("K85").Value = 1703
$ws.Range("L85").Value = 2001.5
$ws.Range("M85").Value = -455
$ws.Range("N85").Value = -4497.5
$ws.Range("H126").Value = 14672.3
$ws.Range("I126").Value = 28800.75
$ws.Range("K126").Value = 86402.25
$ws.Range("M126").Value = -83932.25
$ws.Range("H136").Value = 8378.808000000001
$ws.Range("I136").Value = 6996
$ws.Range("K136").Value = 20988
$ws.Range("M136").Value = -18438

# Worksheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H25").Value = 13899
$ws.Range("J25").Value = 13899
$ws.Range("L25").Value = 13899
$ws.Range("N25").Value = -14485
$ws.Range("H47").Value = 64999.168
$ws.Range("J47").Value = 64999.168
$ws.Range("L47").Value = 64999.168
$ws.Range("N47").Value = -66143.16800000001
$ws.Range("H81").Value = 2997.1333
$ws.Range("I81").Value = 2997.1333
$ws.Range("K81").Value = 5994.2666
$ws.Range("M81").Value = -4933.2666
$ws.Range("H84").Value = 2997.1333
$ws.Range("I84").Value = 2997.1333
$ws.Range("K84").Value = 29971.333
$ws.Range("M84").Value = -24667.333
$ws.Range("H132").Value = 6025.875
$ws.Range("I132").Value = 5427.1665
$ws.Range("K132").Value = 16281.4995
$ws.Range("M132").Value = -13751.4995
